$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.568.83"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "3.024.59"
$ws.Range("E3").Value = "  +2.33%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Formula = "=""380.18"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Formula = "=""102.34"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Formula = "=""0.588"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").Formula = "=""36.61"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").Value = "3.500.22"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("D14").Formula = "=""18.47"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").Value = "3.026.75"
$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").Formula = "=""0.971"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -3.42%  "

$ws.Range("D18").Formula = "=""10.62"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -14.88%  "

$ws.Range("D19").Value = "51.566.97"
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("D20").Formula = "=""3.07"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").Formula = "=""12.43"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").Formula = "=""69.88"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").Formula = "=""267.23"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Formula = "=""3.15"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -5.53%  "

$ws.Range("D26").Formula = "=""8.48"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +6.96%  "

$ws.Range("D27").Formula = "=""7.51"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +6.65%  "

$ws.Range("E28").Value = "  +3.37%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").Formula = "=""26.15"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").Formula = "=""10.25"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").Formula = "=""34.09"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").Formula = "=""50.53"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.17%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Formula = "=""0.0448"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Formula = "=""2.02"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").Formula = "=""0.297"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +12.21%  "

$ws.Range("D40").Formula = "=""16.97"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +1.60%  "

$ws.Range("D41").Formula = "=""129.11"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +3.67%  "

$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("E45").Value = "  +4.48%  "

$ws.Range("D46").Formula = "=""21.58"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("D48").Formula = "=""2.42"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +2.74%  "

$ws.Range("D49").Value = "2.022.14"
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("D50").Value = "3.324.95"
$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("D51").Formula = "=""0.515"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +4.93%  "
